$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51: coin renamed from Quant to Cronos (with updated link)
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"

$ws.Range("D2").Value = "30.807.63"
$ws.Range("E2").Value = "  +2.35%  "
$ws.Range("D3").Value = "2.112.62"
$ws.Range("E3").Value = "  +9.73%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.51"
$ws.Range("E5").Value = "  +4.15%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5284"
$ws.Range("E7").Value = "  +3.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4378"
$ws.Range("E8").Value = "  +8.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08997"
$ws.Range("E9").Value = "  +7.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.03"
$ws.Range("E10").Value = "  +9.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.176"
$ws.Range("E11").Value = "  +5.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.93"
$ws.Range("E12").Value = "  +3.91%  "
$ws.Range("D13").Value = "2.113.42"
$ws.Range("E13").Value = "  +9.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.765"
$ws.Range("E14").Value = "  +5.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.806"
$ws.Range("E15").Value = "  +7.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.23"
$ws.Range("E16").Value = "  +4.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001130"
$ws.Range("E18").Value = "  +2.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06668"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.07"
$ws.Range("E20").Value = "  +3.05%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.362"
$ws.Range("E22").Value = "  +6.88%  "
$ws.Range("D23").Value = "30.885.35"
$ws.Range("E23").Value = "  +2.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.15"
$ws.Range("E24").Value = "  +6.97%  "
$ws.Range("D25").Value = "2.362.40"
$ws.Range("E25").Value = "  +10.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.268"
$ws.Range("E26").Value = "  +3.73%  "
$ws.Range("E27").Value = "  +3.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.564"
$ws.Range("E28").Value = "  +13.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.64"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.87"
$ws.Range("E30").Value = "  +3.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.165"
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.225"
$ws.Range("E33").Value = "  +4.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.055"
$ws.Range("E34").Value = "  +7.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.529"
$ws.Range("E35").Value = "  +20.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02604"
$ws.Range("E36").Value = "  +6.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.534"
$ws.Range("E37").Value = "  +4.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06734"
$ws.Range("E38").Value = "  +4.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.526"
$ws.Range("E39").Value = "  +9.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.77"
$ws.Range("E40").Value = "  +8.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2269"
$ws.Range("E41").Value = "  +5.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6853"
$ws.Range("E42").Value = "  +5.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.251"
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6455"
$ws.Range("E44").Value = "  +6.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "14.06"
$ws.Range("E46").Value = "  +5.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.235"
$ws.Range("E47").Value = "  +3.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.670"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.274"
$ws.Range("E49").Value = "  +5.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.46"
$ws.Range("E50").Value = "  +5.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07112"
$ws.Range("E51").Value = "  +4.12%  "
